$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'" + '26.432.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "'" + '  +2.08%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value2 = "'" + '1.669.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "'" + '  +1.75%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value2 = "'" + '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "'" + '  +0.26%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value2 = "'" + '219.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "'" + '  +3.03%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value2 = "'" + '0.5257'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "'" + '  +0.99%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value2 = "'" + '1.004'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "'" + '  +0.20%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value2 = "'" + '0.2666'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "'" + '  +2.56%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value2 = "'" + '0.06365'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "'" + '  +0.79%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value2 = "'" + '21.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "'" + '  +5.17%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value2 = "'" + '0.07804'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "'" + '  +1.77%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value2 = "'" + 'Polkadot'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value2 = "'" + 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value2 = "'" + '4.463'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "'" + '  +1.22%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value2 = "'" + 'WrappedEther'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value2 = "'" + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value2 = "'" + '1.658.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "'" + '  +1.40%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value2 = "'" + '  +0.94%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value2 = "'" + '  +1.21%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value2 = "'" + '65.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "'" + '  +1.57%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value2 = "'" + '26.442.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "'" + '  +2.07%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value2 = "'" + '  +0.17%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value2 = "'" + '4.737'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "'" + '  +1.09%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value2 = "'" + '193.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "'" + '  +2.85%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value2 = "'" + '10.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "'" + '  +2.03%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value2 = "'" + '6.263'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "'" + '  +0.32%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value2 = "'" + '1.006'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "'" + '  +0.38%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value2 = "'" + 'Monero'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value2 = "'" + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value2 = "'" + '138.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "'" + '  -2.87%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value2 = "'" + 'Stellar'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value2 = "'" + 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value2 = "'" + '0.1259'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "'" + '  +1.78%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value2 = "'" + '7.394'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "'" + '  +0.57%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value2 = "'" + '16.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "'" + '  +2.50%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value2 = "'" + '1.417'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "'" + '  +0.63%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value2 = "'" + '0.06165'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "'" + '  +4.46%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value2 = "'" + '  +2.46%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value2 = "'" + '  +6.49%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value2 = "'" + '3.395'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "'" + '  +0.30%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value2 = "'" + '1.682'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "'" + '  +2.90%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value2 = "'" + '1.003'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "'" + '  +1.68%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value2 = "'" + '0.6080'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "'" + '  +8.81%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value2 = "'" + '2.421'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "'" + '  +1.10%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value2 = "'" + '2.767'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "'" + '  +1.00%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value2 = "'" + '  +0.90%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value2 = "'" + '6.030'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "'" + '  +3.48%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value2 = "'" + '1.090.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "'" + '  +7.05%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value2 = "'" + '0.8587'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "'" + '  +0.97%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value2 = "'" + '  +0.13%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value2 = "'" + '100.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "'" + '  +2.14%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value2 = "'" + '1.812.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "'" + '  +1.31%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value2 = "'" + 'Aave'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value2 = "'" + 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value2 = "'" + '57.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "'" + '  +4.66%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value2 = "'" + 'BabyDogeCoin'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value2 = "'" + 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value2 = "'" + '0.0₈107'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "'" + '  -2.84%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value2 = "'" + '8.176'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "'" + '  +2.27%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value2 = "'" + '1.004'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "'" + '  -0.12%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value2 = "'" + '0.05206'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "'" + '  +1.32%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value2 = "'" + '1.484'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "'" + '  +8.65%  '
$ws.Range("E50").Style = "Normal"
